$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing "SOX" (exact match, whole cell) in column A and
# delete that entire row, shifting the remaining stock symbols up so the
# list stays contiguous (this also removes the now-unused last row).
$found = $ws.UsedRange.Find("SOX", [Type]::Missing, [Type]::Missing, 1)
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
